# Wildlife health database 2-pager — "updated ERM July 11" commit.
#
# Two tracked insertions, both authored by "Montecino, Diego":
#   1. Title: "Data management investments for wildlife health"
#      -> append " in protected areas"
#   2. Body: "...wildlife mortality in Southeast Asia..."
#      -> insert "protected areas of " right before "Southeast Asia"

$word.Application.UserName = "Montecino, Diego"

$d = $word.ActiveDocument
$d.TrackRevisions = $true

# --- Edit 1: title ------------------------------------------------------
$title = $d.Content
$found1 = $title.Find.Execute(
    "Data management investments for wildlife health",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found1) {
    $title.Collapse(0)  # wdCollapseEnd
    $title.InsertAfter(" in protected areas")
}

# --- Edit 2: "Southeast Asia" mention in the disease-burden paragraph --
$ctx = $d.Content
$found2 = $ctx.Find.Execute(
    "wildlife mortality in Southeast Asia",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $insertAt = $ctx.End - 14  # "Southeast Asia" is 14 characters long
    $pt = $d.Range($insertAt, $insertAt)
    $pt.InsertBefore("protected areas of ")
}
